$d = $word.ActiveDocument

# Anchor on the "LOB1036: ..." requirement paragraph, then walk forward
# through the four paragraphs that follow it:
#   1. an empty paragraph
#   2. "Ver no Jupiter Salvar em pdf Salvar em docx"
#   3. an empty paragraph
#   4. an empty paragraph with PageBreakBefore + left alignment
# All four are removed, leaving the trailing empty paragraph and the
# trailing page-break paragraph (plus the sectPr) untouched.

$anchor = $d.Content
$found = $anchor.Find.Execute("LOB1036: Geometria Analítica (Requisito fraco)", `
    $false, $false, $false, $false, $false, $true, 1, $false, "", 0)

$reqPara = $anchor.Paragraphs(1)

$p1 = $reqPara.Next()       # empty paragraph
$p2 = $p1.Next()            # "Ver no Jupiter Salvar em pdf Salvar em docx"
$p3 = $p2.Next()            # empty paragraph
$p4 = $p3.Next()            # empty paragraph, PageBreakBefore + jc=left

$delRange = $d.Range($p1.Range.Start, $p4.Range.End)
$delRange.Delete()
